$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 34, column B: change "Yes" -> "No"
$ws.Cells.Item(34, 2).Value = "No"

# Copy formatting from row 34 into the new row 35 so styles match
$src = $ws.Range("A34:I34")
$dst = $ws.Range("A35:I35")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new row 35 values
$ws.Cells.Item(35, 1).Value = "Verify Add Credit Card With Invalid Data"
$ws.Cells.Item(35, 2).Value = "Yes"
$ws.Cells.Item(35, 3).Value = "testdata.xls,PaymentMethods"
$ws.Cells.Item(35, 4).Value = "RunOneIteration"
$ws.Cells.Item(35, 5).Value = "'1"
$ws.Cells.Item(35, 6).Value = "'1"
$ws.Cells.Item(35, 7).Value = "PaymentMethods"
$ws.Cells.Item(35, 8).Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin,`n-puserName"
$ws.Cells.Item(35, 9).Value = "coyni_mobile.tests.CustomerProfileTest,`ntestAddDebitCardWithInvalidData,`n-pnameOnCard,`n-pcardNumber,`n-pcardExp,`n-pcvvOrCVC,`n-paddressLine1,`n-paddreddLine2,`n-pcity,`n-pstate,`n-pzipCode,`n-perrMsg,`n-pelementName"

# Match the autofit row height Excel would compute for the 13-line wrapped text
$ws.Rows.Item(35).RowHeight = 187.2

# Update sheet view to reflect scrolled/selected state for the new row
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("I35").Select()
